$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 -> becomes what was row 25 (for columns A,B,E,F,G,H,Q,R)
$ws.Range("A24").Value = 130827882
$ws.Range("B24").Value = 8440
$ws.Range("E24").Value = 106554
$ws.Range("F24").Value = "Björksplintborre"
$ws.Range("G24").Value = "Scolytus ratzeburgii"
$ws.Range("H24").Value = "Janson, 1856"
$ws.Range("Q24").Value = 344503
$ws.Range("R24").Value = 6433291

# Row 24's AJ/AK/AO (gran / Picea abies / Picea abies) are removed (moved to row 25)
$ws.Range("AJ24").ClearContents()
$ws.Range("AK24").ClearContents()
$ws.Range("AO24").ClearContents()

# Row 25 -> becomes what was row 24 (for columns A,B,E,F,G,H,Q,R)
$ws.Range("A25").Value = 130827870
$ws.Range("B25").Value = 5197
$ws.Range("E25").Value = 105930
$ws.Range("F25").Value = "Vågbandad barkbock"
$ws.Range("G25").Value = "Semanotus undatus"
$ws.Range("H25").Value = "(Linnaeus, 1758)"
$ws.Range("Q25").Value = 344458
$ws.Range("R25").Value = 6433350

# Row 25 gains AJ/AK/AO values that used to belong to row 24
$ws.Range("AJ25").Value = "gran"
$ws.Range("AK25").Value = "Picea abies"
$ws.Range("AO25").Value = "Picea abies"
